# Update country data figures on the "Canada Summary" worksheet.
#
# In the source workbook these figures are stored as plain text (shared
# strings), not numbers -- e.g. "34.7", "29", "99.9" -- with General
# number formatting. We keep that text representation by prefixing the
# new values with a leading apostrophe, which tells Excel to treat the
# entry as text even though it looks numeric, instead of silently
# converting the cell to a numeric type.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enterprises density (per 1000 people) row (row 13): Micro/SMEs/MSMEs
$ws.Range("B13").Value = "'34.74"
$ws.Range("C13").Value = "'14.21"
$ws.Range("D13").Value = "'48.95"

# Employment (% of total) row (row 14): Micro/MSMEs
$ws.Range("B14").Value = "'5.13"
$ws.Range("D14").Value = "'90.43"

# Enterprises (% of total) row (row 16): SMEs/MSMEs
$ws.Range("C16").Value = "'29.01"
$ws.Range("D16").Value = "'99.91"
